$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing row, values updated)
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr10"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.380435
$ws.Range("H2").Value = 1.141305
$ws.Range("I2").Value = 0.7997108917301441
$ws.Range("J2").Value = 0.7997108917301442
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.036595333333333
$ws.Range("N2").Value = 3.109786
$ws.Range("O2").Value = 0.393072250513715
$ws.Range("P2").Value = 0.393072250513715
$ws.Range("Q2").Value = 0.3943571456366667
$ws.Range("R2").Value = 3.54921431073
$ws.Range("S2").Value = 0.3143441599726977
$ws.Range("T2").Value = 0.3143441599726977

# Row 3 (new row)
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.380435
$ws.Range("H3").Value = 1.141305
$ws.Range("I3").Value = 0.7997108917301441
$ws.Range("J3").Value = 0.7997108917301442
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.600567
$ws.Range("N3").Value = 4.801701
$ws.Range("O3").Value = 0.6069277494862849
$ws.Range("P3").Value = 0.6069277494862849
$ws.Range("Q3").Value = 0.6089117066450001
$ws.Range("R3").Value = 5.480205359805001
$ws.Range("S3").Value = 0.4853667317574464
$ws.Range("T3").Value = 0.4853667317574465

# Row 4 (new row)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.09528066666666667
$ws.Range("H4").Value = 0.285842
$ws.Range("I4").Value = 0.2002891082698559
$ws.Range("J4").Value = 0.2002891082698559
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.036595333333333
$ws.Range("N4").Value = 3.109786
$ws.Range("O4").Value = 0.393072250513715
$ws.Range("P4").Value = 0.393072250513715
$ws.Range("Q4").Value = 0.09876749442355555
$ws.Range("R4").Value = 0.888907449812
$ws.Range("S4").Value = 0.07872809054101737
$ws.Range("T4").Value = 0.07872809054101738

# Row 5 (new row)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.09528066666666667
$ws.Range("H5").Value = 0.285842
$ws.Range("I5").Value = 0.2002891082698559
$ws.Range("J5").Value = 0.2002891082698559
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.600567
$ws.Range("N5").Value = 4.801701
$ws.Range("O5").Value = 0.6069277494862849
$ws.Range("P5").Value = 0.6069277494862849
$ws.Range("Q5").Value = 0.1525030908046667
$ws.Range("R5").Value = 1.372527817242
$ws.Range("S5").Value = 0.1215610177288385
$ws.Range("T5").Value = 0.1215610177288385
